# Apply the "break out stock.yaml completed" edit to the "day" sheet:
#  1. Convert D77:D84 (bsecode) from text to real numbers.
#  2. Append 5 new rows (85-89) of stock data pulled from the 04/07/2024 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day")

# --- 1. Fix up D77:D84 so the bsecode column holds numbers, not text ---
$bsecodes = @{
    77 = 541154
    78 = 502355
    79 = 500420
    80 = 500087
    81 = 532187
    82 = 500575
    83 = 511243
    84 = 540777
}
foreach ($row in $bsecodes.Keys) {
    $ws.Cells.Item($row, 4).Value = $bsecodes[$row]
}

# --- 2. Append the new rows (85-89) ---
# Note: the bsecode column (D) stays text here (a leading "'" keeps the
# numeric-looking string from being auto-converted to a number).
$newRows = @(
    @(1, "COROMANDEL", "Coromandel International Limited", "'506395", -0.26, 1579.55, 154337, "day", "04/07/2024 11:34:41"),
    @(2, "GNFC", "Gujarat Narmada Valley Fertilizers And Chemicals Limited", "'500670", 0.62, 727.95, 2335100, "day", "04/07/2024 11:34:41"),
    @(3, "VEDL", "Vedanta Limited", "'500295", 1.12, 469.1, 12818706, "day", "04/07/2024 11:34:41"),
    @(4, "NMDC", "Nmdc Limited", "'526371", 0.16, 251.61, 13292696, "day", "04/07/2024 11:34:41"),
    @(5, "ABCAPITAL", "Aditya Birla Capital Ltd", "'540691", -0.08, 237.38, 3314009, "day", "04/07/2024 11:34:41")
)

$startRow = 85
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
    $ws.Cells.Item($r, 7).Value = $data[6]
    $ws.Cells.Item($r, 8).Value = $data[7]
    $ws.Cells.Item($r, 9).Value = $data[8]
    # Undo the "quote prefix" style Excel applies when a value is entered
    # with a leading apostrophe, so the cell keeps the default look.
    $ws.Cells.Item($r, 4).Style = "Normal"
}
